$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the four "raw" category labels in row 58 to their friendlier
#    chart-facing names. (B/C/E/F/H/I keep their existing text; only the
#    L/M/N/O labels actually change text.)
# ---------------------------------------------------------------------------
$ws.Range("L58").Value = "Errors on Ribbon"
$ws.Range("M58").Value = "Errors on CommandMaps"
$ws.Range("N58").Value = "Rating for Ribbon"
$ws.Range("O58").Value = "Rating for CommandMaps"

# ---------------------------------------------------------------------------
# 2. Fix the standard-error formulas on row 61: divide by sqrt(53) (the
#    correct post-filter sample size) instead of sqrt(71), for every column
#    except D and G (whose source stddev, row 60, is 0 and which keep
#    referencing the original formula).
# ---------------------------------------------------------------------------
foreach ($col in @("B","C","E","F","H","I","L","M","N","O")) {
    $ws.Range($col + "61").Formula = "=" + $col + "60/(53 ^ 0.5)"
}

# ---------------------------------------------------------------------------
# 3. Add titles to the three charts.
# ---------------------------------------------------------------------------
$chartTime = $ws.ChartObjects(1).Chart
$chartTime.HasTitle = $true
$chartTime.ChartTitle.Text = "Mean Time per Trial"

$chartErrors = $ws.ChartObjects(2).Chart
$chartErrors.HasTitle = $true
$chartErrors.ChartTitle.Text = "Mean Errors after 72 Trial Experiment"

$chartRating = $ws.ChartObjects(3).Chart
$chartRating.HasTitle = $true
$chartRating.ChartTitle.Text = "Mean Rating"

# ---------------------------------------------------------------------------
# 4. Resize/reposition the Errors and Rating charts to account for the new
#    titles (matches the taller chart frames saved by Excel after a title
#    is added).
# ---------------------------------------------------------------------------
$errorsChartObj = $ws.ChartObjects(2)
$errorsChartObj.Left = 72.4375
$errorsChartObj.Top = 975.5
$errorsChartObj.Width = 342.625
$errorsChartObj.Height = 246.5

$ratingChartObj = $ws.ChartObjects(3)
$ratingChartObj.Left = 71.4375
$ratingChartObj.Top = 1237.5
$ratingChartObj.Width = 347.625
$ratingChartObj.Height = 230.5

# ---------------------------------------------------------------------------
# 5. Restore the author's final selection/scroll state.
# ---------------------------------------------------------------------------
$null = $ws.Range("T73").Select()
